# Apply cryptocurrency price/volume update (GitHub Actions style refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking Price cells to stay text (matches the source data,
# which stores prices like "0.4576" as text, not as a number) before writing them.
# (Applied per contiguous block, since a comma-unioned Range only formats its first area.)
$ws.Range("D4:D12").NumberFormat = "@"
$ws.Range("D14:D16").NumberFormat = "@"
$ws.Range("D18:D19").NumberFormat = "@"
$ws.Range("D23:D24").NumberFormat = "@"
$ws.Range("D26:D34").NumberFormat = "@"
$ws.Range("D36:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.810.68"
$ws.Range("E2").Value = "  -1.84%  "
$ws.Range("D3").Value = "1.890.94"
$ws.Range("E3").Value = "  -2.94%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").Value = "323.37"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("D7").Value = "0.4576"
$ws.Range("E7").Value = "  -1.16%  "
$ws.Range("D8").Value = "0.3802"
$ws.Range("E8").Value = "  -2.98%  "
$ws.Range("D9").Value = "45.38"
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("D10").Value = "0.07714"
$ws.Range("E10").Value = "  -2.30%  "
$ws.Range("D11").Value = "0.9623"
$ws.Range("E11").Value = "  -3.80%  "
$ws.Range("D12").Value = "21.95"
$ws.Range("E12").Value = "  -1.96%  "
$ws.Range("D13").Value = "1.876.90"
$ws.Range("E13").Value = "  -3.60%  "
$ws.Range("D14").Value = "6.971"
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").Value = "5.668"
$ws.Range("E15").Value = "  -3.29%  "
$ws.Range("D16").Value = "0.07067"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "83.40"
$ws.Range("E18").Value = "  -5.33%  "
$ws.Range("D19").Value = "0.000009504"
$ws.Range("E19").Value = "  -4.76%  "
$ws.Range("E20").Value = "  -2.16%  "
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").Value = "28.795.10"
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("D23").Value = "5.367"
$ws.Range("E23").Value = "  -2.94%  "
$ws.Range("D24").Value = "10.92"
$ws.Range("E24").Value = "  -2.67%  "
$ws.Range("D25").Value = "2.116.48"
$ws.Range("E25").Value = "  -2.93%  "
$ws.Range("D26").Value = "2.080"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").Value = "155.64"
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").Value = "19.08"
$ws.Range("E28").Value = "  -2.32%  "
$ws.Range("D29").Value = "5.647"
$ws.Range("E29").Value = "  -4.92%  "
$ws.Range("D30").Value = "117.18"
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("D31").Value = "1.817"
$ws.Range("E31").Value = "  -4.03%  "
$ws.Range("D32").Value = "0.09259"
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("D33").Value = "0.8533"
$ws.Range("E33").Value = "  -5.17%  "
$ws.Range("D34").Value = "5.070"
$ws.Range("E34").Value = "  -3.19%  "
$ws.Range("E35").Value = "  -6.57%  "
$ws.Range("D36").Value = "3.071"
$ws.Range("E36").Value = "  -2.94%  "
$ws.Range("D37").Value = "1.154"
$ws.Range("E37").Value = "  -1.89%  "
$ws.Range("D38").Value = "0.05644"
$ws.Range("E38").Value = "  -2.79%  "
$ws.Range("D39").Value = "1.005"
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("D40").Value = "0.02042"
$ws.Range("E40").Value = "  -3.48%  "
$ws.Range("D41").Value = "0.5512"
$ws.Range("E41").Value = "  -3.92%  "
$ws.Range("D42").Value = "7.431"
$ws.Range("E42").Value = "  -3.82%  "
$ws.Range("D43").Value = "0.1754"
$ws.Range("E43").Value = "  -3.68%  "
$ws.Range("D44").Value = "0.000002885"
$ws.Range("E44").Value = "  -25.18%  "
$ws.Range("D45").Value = "9.253"
$ws.Range("E45").Value = "  -5.66%  "
$ws.Range("D46").Value = "2.696"
$ws.Range("E46").Value = "  +3.64%  "
$ws.Range("D47").Value = "0.5172"
$ws.Range("E47").Value = "  -3.46%  "
$ws.Range("D48").Value = "11.17"
$ws.Range("E48").Value = "  -6.71%  "
$ws.Range("D49").Value = "2.062"
$ws.Range("E49").Value = "  -7.19%  "
$ws.Range("D50").Value = "0.06763"
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("D51").Value = "110.63"
$ws.Range("E51").Value = "  -2.83%  "
